$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.535.59'
$ws.Range('E2').Value = '  +6.14%  '
$ws.Range('D3').Value = '3.486.58'
$ws.Range('E3').Value = '  +7.48%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.00'
$ws.Range('E5').Value = '  +7.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.24'
$ws.Range('E6').Value = '  +8.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '3.495.22'
$ws.Range('E8').Value = '  +7.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.62'
$ws.Range('E10').Value = '  +3.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.124'
$ws.Range('E11').Value = '  +8.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.444'
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('D13').Value = '4.079.10'
$ws.Range('E13').Value = '  +7.22%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000189'
$ws.Range('E15').Value = '  +8.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.85'
$ws.Range('E16').Value = '  +5.85%  '
$ws.Range('D17').Value = '64.579.92'
$ws.Range('E17').Value = '  +6.25%  '
$ws.Range('D18').Value = '3.487.96'
$ws.Range('E18').Value = '  +7.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.49'
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.47'
$ws.Range('E20').Value = '  +7.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '399.76'
$ws.Range('E21').Value = '  +6.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.56'
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.549'
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.43'
$ws.Range('E25').Value = '  +3.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000112'
$ws.Range('E26').Value = '  +22.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.49'
$ws.Range('E27').Value = '  +10.12%  '
$ws.Range('E28').Value = '  +6.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.42'
$ws.Range('E30').Value = '  +15.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.98'
$ws.Range('E31').Value = '  +10.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.76'
$ws.Range('E32').Value = '  +9.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.06'
$ws.Range('E33').Value = '  +6.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.96'
$ws.Range('E34').Value = '  +6.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.98'
$ws.Range('E36').Value = '  +5.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.52'
$ws.Range('E37').Value = '  +6.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.80'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '28.75'
$ws.Range('E39').Value = '  +8.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0791'
$ws.Range('E40').Value = '  +9.86%  '
$ws.Range('E41').Value = '  +9.99%  '
$ws.Range('D42').Value = '2.892.78'
$ws.Range('E42').Value = '  +3.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0325'
$ws.Range('E43').Value = '  +3.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.786'
$ws.Range('E44').Value = '  +7.72%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.45'
$ws.Range('E45').Value = '  +3.95%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.29'
$ws.Range('E46').Value = '  +5.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  +11.08%  '
$ws.Range('D48').Value = '3.530.78'
$ws.Range('E48').Value = '  +7.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.91'
$ws.Range('E49').Value = '  +6.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '300.60'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.13'
$ws.Range('E51').Value = '  +22.97%  '
